$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Coin) and C (Link) text replacements, rows 8-17 ---
$ws.Range("B8").Value = "MXToken"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("B10").Value = "WazirX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("B16").Value = "LEO"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"

# --- Columns D (Price) and E (Volume/1h) updates ---
# These new values look like numbers/percentages, so Excel would normally convert
# them into numeric cells. The source file stores them as plain text, so we force
# text storage via NumberFormat "@" and then restore the default "Normal" style
# afterwards, leaving the cell visually/structurally as it was (no explicit style).
$fmtRange = $ws.Range("D2:D27")
$fmtRange.NumberFormat = "@"
$ws.Range("D2").Value = "307.69"
$ws.Range("D3").Value = "40.02"
$ws.Range("D4").Value = "5.129"
$ws.Range("D5").Value = "0.08099"
$ws.Range("D6").Value = "1.939"
$ws.Range("D7").Value = "8.121"
$ws.Range("D8").Value = "0.9302"
$ws.Range("D9").Value = "0.1436"
$ws.Range("D10").Value = "0.1913"
$ws.Range("D11").Value = "0.09053"
$ws.Range("D12").Value = "0.03530"
$ws.Range("D13").Value = "0.09809"
$ws.Range("D14").Value = "0.001404"
$ws.Range("D15").Value = "0.005904"
$ws.Range("D16").Value = "3.888"
$ws.Range("D17").Value = "4.221"
$ws.Range("D18").Value = "3.393"
$ws.Range("D19").Value = "0.3424"
$ws.Range("D20").Value = "0.1320"
$ws.Range("D21").Value = "4.667"
$ws.Range("D22").Value = "0.2500"
$ws.Range("D23").Value = "0.04390"
$ws.Range("D24").Value = "0.001217"
$ws.Range("D25").Value = "0.004332"
$ws.Range("D26").Value = "0.0001305"
$ws.Range("D27").Value = "0.0004013"
$fmtRange.Style = "Normal"

$fmtRange = $ws.Range("D39:D48")
$fmtRange.NumberFormat = "@"
$ws.Range("D39").Value = "0.02031"
$ws.Range("D40").Value = "0.05038"
$ws.Range("D41").Value = "0.007399"
$ws.Range("D42").Value = "0.009916"
$ws.Range("D43").Value = "0.1363"
$ws.Range("D44").Value = "0.002138"
$ws.Range("D45").Value = "0.009197"
$ws.Range("D46").Value = "0.00006368"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("D48").Value = "0.002873"
$fmtRange.Style = "Normal"

$fmtRange = $ws.Range("D50:D51")
$fmtRange.NumberFormat = "@"
$ws.Range("D50").Value = "0.00002107"
$ws.Range("D51").Value = "0.0002006"
$fmtRange.Style = "Normal"

$fmtRange = $ws.Range("E2:E27")
$fmtRange.NumberFormat = "@"
$ws.Range("E2").Value = "-0.03%"
$ws.Range("E3").Value = "5.42%"
$ws.Range("E4").Value = "1.03%"
$ws.Range("E5").Value = "-0.58%"
$ws.Range("E6").Value = "-1.86%"
$ws.Range("E7").Value = "2.83%"
$ws.Range("E8").Value = "0.17%"
$ws.Range("E9").Value = "1.56%"
$ws.Range("E10").Value = "-1.65%"
$ws.Range("E11").Value = "-2.27%"
$ws.Range("E12").Value = "0.81%"
$ws.Range("E13").Value = "-0.82%"
$ws.Range("E14").Value = "-0.06%"
$ws.Range("E15").Value = "-4.53%"
$ws.Range("E16").Value = "-1.46%"
$ws.Range("E17").Value = "1.14%"
$ws.Range("E18").Value = "-0.54%"
$ws.Range("E19").Value = "-0.87%"
$ws.Range("E20").Value = "2.72%"
$ws.Range("E21").Value = "-3.00%"
$ws.Range("E22").Value = "-4.51%"
$ws.Range("E23").Value = "-2.23%"
$ws.Range("E24").Value = "-2.14%"
$ws.Range("E25").Value = "3.80%"
$ws.Range("E26").Value = "0.25%"
$ws.Range("E27").Value = "-9.78%"
$fmtRange.Style = "Normal"

$fmtRange = $ws.Range("E39:E47")
$fmtRange.NumberFormat = "@"
$ws.Range("E39").Value = "-4.28%"
$ws.Range("E40").Value = "-2.17%"
$ws.Range("E41").Value = "-0.97%"
$ws.Range("E42").Value = "-2.26%"
$ws.Range("E43").Value = "-0.42%"
$ws.Range("E44").Value = "0.26%"
$ws.Range("E45").Value = "-5.09%"
$ws.Range("E46").Value = "0.65%"
$ws.Range("E47").Value = "0.21%"
$fmtRange.Style = "Normal"

$fmtRange = $ws.Range("E49:E51")
$fmtRange.NumberFormat = "@"
$ws.Range("E49").Value = "-18.81%"
$ws.Range("E50").Value = "0.21%"
$ws.Range("E51").Value = "0.21%"
$fmtRange.Style = "Normal"

# --- Column G (Hora) updates: every changed row becomes "3" ---
$fmtRange = $ws.Range("G2:G51")
$fmtRange.NumberFormat = "@"
$fmtRange.Value = "3"
$fmtRange.Style = "Normal"
